# route.xlsx template update
#  - Pass client timezone to excel template (use joda DateTime with timezone
#    for position.fixTime, and a plain toString(...) concatenation for the
#    from/to period instead of the old String.format(...) call)
#  - Templates: adjust timezone / formatting, use "https" in links
#  - Minor layout tweaks: paragraph indents, column widths, active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Period: ${"".format(...)}  ->  ${from.toString(...)+" - "+to.toString(...)}
$ws.Range("B6").Formula = '${from.toString("YYYY.MM.dd HH:mm:ss")+" - "+to.toString("YYYY.MM.dd HH:mm:ss")}'

# --- Time column body cell: ${position.fixTime} -> joda DateTime with timezone
$ws.Range("B9").Formula = '${new("org.joda.time.DateTime", position.fixTime, timezone).toString("YYYY.MM.dd HH:mm:ss")}'

# --- Address hyperlink: http -> https
$ws.Range("G9").Formula = '${util.hyperlink("".format("https://www.openstreetmap.org/?mlat=%1$f&mlon=%2$f#map=16/%1$f/%2$f", position.latitude, position.longitude), position.getAddress() == null ? "".format("%1$f°, %2$f°", position.latitude, position.longitude) : position.address)}'

# --- Attributes column (unchanged content, re-asserted for safety)
$ws.Range("H9").Formula = '${position.attributes.toString().replaceAll(",", " ").replaceAll(bracketsRegex, "")}'

# --- Paragraph indents: bump a few label cells over
$ws.Range("B1").IndentLevel = 15
$ws.Range("B3").IndentLevel = 15
$ws.Range("B2").IndentLevel = 2
$ws.Range("B4").IndentLevel = 2
$ws.Range("B5").IndentLevel = 2
$ws.Range("B6").IndentLevel = 2

# --- Column widths (G/H slightly wider)
$ws.Range("G1").EntireColumn.ColumnWidth = 61.0
$ws.Range("H1").EntireColumn.ColumnWidth = 72.5

# --- Active selection moves to G9
$ws.Range("G9").Select()

# --- Tab ratio (sheet-tabs / horizontal-scrollbar split) nudged slightly
$win = $excel.ActiveWindow
$win.TabRatio = 0.989
